$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 572, shifting existing rows 572-639 down to 573-640.
$ws.Rows.Item(572).Insert()

# Populate the newly inserted row 572 with the new weekly data point.
$ws.Range("A572").Value = 3
$ws.Range("B572").Value = "Femacal de La Calera"
$ws.Range("C572").Value = "Coquimbo"
$ws.Range("D572").Value = 45212
$ws.Range("D572").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E572").Value = 5
$ws.Range("F572").Value = 100114013
$ws.Range("G572").Value = "Zanahoria"
$ws.Range("H572").Value = "Sin especificar"
$ws.Range("I572").Value = "Primera"
$ws.Range("J572").Value = 280
$ws.Range("K572").Value = 6500
$ws.Range("L572").Value = 7000
$ws.Range("M572").Value = 6786
$ws.Range("N572").Value = "$/saco 20 kilos"
$ws.Range("O572").Value = "Provincia de Quillota"
$ws.Range("P572").Value = 339
$ws.Range("Q572").Value = 20
$ws.Range("R572").Value = "Hortaliza"
